$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi Comp -> Cd47 results recomputed following Dr Hou's advice: the
# sending/target cluster pairing now spans all three clusters (ECs, FAPs,
# sCs) -- 9 rows total (rows 2-10) -- with refreshed statistics for every
# pair. Columns A-T follow the existing header layout:
#   A Sending cluster          K Receptor-expressing cells
#   B Ligand symbol            L Receptor detection rate
#   C Receptor symbol          M Receptor average expression value
#   D Target cluster           N Receptor total expression value
#   E Ligand-expressing cells  O Receptor derived specificity (avg)
#   F Ligand detection rate    P Receptor derived specificity (total)
#   G Ligand average expr.     Q Edge average expression weight
#   H Ligand total expr.       R Edge total expression weight
#   I Ligand derived spec.(avg)  S Edge average expression derived specificity
#   J Ligand derived spec.(tot)  T Edge total expression derived specificity
$data = @(
    @("ECs", "Comp", "Cd47", "ECs", 2, 0.6666666666666666, 0.4181183333333333, 1.254355, 0.01571144052599341, 0.01571144052599341, 3, 1, 61.27353666666666, 183.82061, 0.3474604587406809, 0.3474604587406808, 25.61958902850555, 230.57630125655, 0.005459104332638594, 0.005459104332638592),
    @("ECs", "Comp", "Cd47", "FAPs", 2, 0.6666666666666666, 0.4181183333333333, 1.254355, 0.01571144052599341, 0.01571144052599341, 3, 1, 89.66709900000001, 269.001297, 0.5084702637939138, 0.5084702637939138, 37.491457988715, 337.423121898435, 0.007988800308834256, 0.007988800308834254),
    @("ECs", "Comp", "Cd47", "sCs", 2, 0.6666666666666666, 0.4181183333333333, 1.254355, 0.01571144052599341, 0.01571144052599341, 3, 1, 25.40615466666667, 76.218464, 0.1440692774654054, 0.1440692774654053, 10.62277904563556, 95.60501141071998, 0.002263535884520559, 0.002263535884520558),
    @("FAPs", "Comp", "Cd47", "ECs", 3, 1, 24.89087966666667, 74.672639, 0.9353131502385497, 0.9353131502385496, 3, 1, 61.27353666666666, 183.82061, 0.3474604587406809, 0.3474604587406808, 1525.152227921088, 13726.37005128979, 0.3249843362480778, 0.3249843362480778),
    @("FAPs", "Comp", "Cd47", "FAPs", 3, 1, 24.89087966666667, 74.672639, 0.9353131502385497, 0.9353131502385496, 3, 1, 89.66709900000001, 269.001297, 0.5084702637939138, 0.5084702637939138, 2231.892971268087, 20087.03674141279, 0.4755789242317119, 0.4755789242317118),
    @("FAPs", "Comp", "Cd47", "sCs", 3, 1, 24.89087966666667, 74.672639, 0.9353131502385497, 0.9353131502385496, 3, 1, 25.40615466666667, 76.218464, 0.1440692774654054, 0.1440692774654053, 632.3815386007218, 5691.433847406496, 0.13474988975876, 0.13474988975876),
    @("sCs", "Comp", "Cd47", "ECs", 3, 1, 1.303350666666667, 3.910052, 0.04897540923545694, 0.04897540923545693, 3, 1, 61.27353666666666, 183.82061, 0.3474604587406809, 0.3474604587406808, 79.86090486352445, 718.74814377172, 0.01701701815996445, 0.01701701815996444),
    @("sCs", "Comp", "Cd47", "FAPs", 3, 1, 1.303350666666667, 3.910052, 0.04897540923545694, 0.04897540923545693, 3, 1, 89.66709900000001, 269.001297, 0.5084702637939138, 0.5084702637939138, 116.867673259716, 1051.809059337444, 0.02490253925336767, 0.02490253925336767),
    @("sCs", "Comp", "Cd47", "sCs", 3, 1, 1.303350666666667, 3.910052, 0.04897540923545694, 0.04897540923545693, 3, 1, 25.40615466666667, 76.218464, 0.1440692774654054, 0.1440692774654053, 33.11312862223645, 298.018157600128, 0.007055851822124822, 0.00705585182212482)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
